$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 75
$ws.Range("D4").Value = 73.2
$ws.Range("E4").Value = 76.7
$ws.Range("C10").Value = 42.2
$ws.Range("D10").Value = 43.3
$ws.Range("E10").Value = 41.3
$ws.Range("C13").Value = 77.955555555555605
$ws.Range("D13").Value = 80.816666666666706
$ws.Range("E13").Value = 75.133333333333397
$ws.Range("D38").Value = 82.290000000000106
$ws.Range("C62").Value = 82.630769230769303
$ws.Range("D62").Value = 85.098717948717905
$ws.Range("E62").Value = 80.394871794871804
$ws.Range("C63").Value = 40.226751592356699
$ws.Range("D63").Value = 39.480254777070101
$ws.Range("E63").Value = 40.6732484076433
$ws.Range("C64").Value = 57.809375000000003
$ws.Range("D64").Value = 55.381250000000001
$ws.Range("E64").Value = 59.578125
$ws.Range("C65").Value = 73.073684210526295
$ws.Range("D65").Value = 72.592105263157904
$ws.Range("E65").Value = 72.842105263157904
$ws.Range("I65").Value = 85.936363636363694
$ws.Range("C66").Value = 54.301276595744703
$ws.Range("D66").Value = 54.621702127659603
$ws.Range("E66").Value = 53.857446808510701
$ws.Range("H66").Value = 53.465656565656502
$ws.Range("D67").Value = 80.694117647058803
$ws.Range("I67").Value = 88.728571428571399
$ws.Range("D68").Value = 86.8333333333333
$ws.Range("H69").Value = 72.683333333333294
$ws.Range("C73").Value = 78.892857142857196
$ws.Range("D73").Value = 80.967857142857198
$ws.Range("E73").Value = 76.807142857142907
$ws.Range("C76").Value = 71.599999999999994
$ws.Range("D76").Value = 72.825000000000003
$ws.Range("E76").Value = 70.733333333333405
$ws.Range("C77").Value = 55.193333333333399
$ws.Range("D77").Value = 53.62
$ws.Range("E77").Value = 56.1933333333333
$ws.Range("C78").Value = 8.0346153846153801
$ws.Range("D78").Value = 7.2923076923076904
$ws.Range("C79").Value = 12.5129032258065
$ws.Range("D79").Value = 12.3032258064516
$ws.Range("E79").Value = 12.8403225806452
$ws.Range("I80").Value = 77.433333333333294
$ws.Range("C81").Value = 59.1933333333333
$ws.Range("D81").Value = 55.626666666666701
$ws.Range("E81").Value = 60.566666666666698
$ws.Range("I81").Value = 90.533333333333303
$ws.Range("C82").Value = 82.348571428571503
$ws.Range("D82").Value = 84.821428571428598
$ws.Range("E82").Value = 80.099999999999994
$ws.Range("C83").Value = 38.2232394366197
$ws.Range("D83").Value = 37.774647887324001
$ws.Range("E83").Value = 38.571830985915497
$ws.Range("D84").Value = 93.657894736842195
$ws.Range("C87").Value = 75.174285714285702
$ws.Range("D87").Value = 75.7628571428572
$ws.Range("E87").Value = 74.948571428571398
$ws.Range("C88").Value = 48.1666666666667
$ws.Range("D88").Value = 47.4166666666667
$ws.Range("E88").Value = 48.466666666666697
$ws.Range("C89").Value = 50.1413043478261
$ws.Range("D89").Value = 47.5717391304348
$ws.Range("E89").Value = 51.05
$ws.Range("I89").Value = 93.200000000000102
$ws.Range("C90").Value = 13.34
$ws.Range("D90").Value = 12.7314285714286
$ws.Range("E90").Value = 13.8828571428571
$ws.Range("D91").Value = 92.196296296296396
$ws.Range("E91").Value = 86.044444444444395
$ws.Range("H91").Value = 69.120833333333294
$ws.Range("C94").Value = 59.9304347826087
$ws.Range("D94").Value = 57.056521739130403
$ws.Range("E94").Value = 61.9652173913044
$ws.Range("C95").Value = 86.274074074074093
$ws.Range("D95").Value = 89.118518518518499
$ws.Range("E95").Value = 84.085185185185196
$ws.Range("I95").Value = 85.183333333333294
$ws.Range("C96").Value = 66.099999999999994
$ws.Range("D96").Value = 67.929411764705904
$ws.Range("E96").Value = 65.276470588235298
$ws.Range("D97").Value = 91.372413793103405
$ws.Range("C98").Value = 80.895652173913007
$ws.Range("D98").Value = 81.547826086956505
$ws.Range("E98").Value = 80.273913043478302
